# Update formulas for GermSpeedCorrected, WeightGermPercent, MeanGermPercent
# and MeanGermNumber on Sheet1 (column C holds the "Details" / formula text).
#
# Commit message: "Updated formulas for MeanGermPercent, MeanGermNumber and
# WeightGermPercent" (GermSpeedCorrected picked up a small wording fix too).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$nl = [char]10

# Row 19 - GermSpeedCorrected: "$FGP$:" -> "$FGP$ is"
$ws.Range("C19").Value = "It is computed as follows:$nl" + `
    "`$`$S_{corrected} = \frac{S}{FGP}`$`$$nl" + `
    "Where, `$FGP`$ is the final germination percentage or germinability."

# Row 21 - MeanGermPercent: G -> GP, T_n -> k
$ws.Range("C21").Value = "It is estimated as follows:$nl" + `
    "`$`$\overline{GP} = \frac{GP}{k}`$`$$nl" + `
    "Where, `$GP`$ is the final germination percentage and `$k`$ is the total number of time intervals (e.g. days) required for final germination."

# Row 22 - MeanGermNumber: T_n -> k, clarified wording
$ws.Range("C22").Value = "It is estimated as follows:$nl" + `
    "`$`$\overline{N} = \frac{N_{g}}{k}`$`$$nl" + `
    "Where, `$N_{g}`$ is the number of germinated seeds at the end of the germination test and `$k`$ is the total number of time intervals (e.g. days) required for final germination."

# Row 20 - WeightGermPercent: interval index $t$ -> $k$
$ws.Range("C20").Value = "It is estimated as follows:$nl" + `
    "`$`$WGP = \frac{\sum_{i=1}^{k}(k-i+1)N_{i}}{k \times N} \times 100`$`$$nl" + `
    "Where, `$N_{i}`$ is the number of seeds that germinated in the time interval `$i`$ (not cumulative, but partial count), `$N`$ is the total number of seeds tested and `$k`$ is the total number of time intervals."

# Reflect the cursor/viewport position captured in the saved file.
$ws.Activate()
$ws.Range("I21").Select()
